$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Patente / Motor / Chasis values in row 2 with new data
$ws.Range("F2").Value = "ZZZ115"
$ws.Range("G2").Value = "ZZZ111BB013"
$ws.Range("H2").Value = "ZZZ111BB1231"

# Remove the now-duplicate row 3
$ws.Rows.Item(3).Delete()

# Update selection to reflect the edited range
$ws.Range("F2:H2").Select()
